$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells hold numeric-looking values but are stored as TEXT (shared
# strings), matching the source workbook. A leading apostrophe forces text
# entry; ClearFormats() strips the quote-prefix cell style Excel applies
# automatically so the cell format stays identical to the original ("General").

# G2V row (row 3): Trading rev. changes (Adj. rev. mirrors the same figure)
$ws.Range("C3").Value = "'-6.5"
$ws.Range("C3").ClearFormats()

$ws.Range("F3").Value = "'-6.5"
$ws.Range("F3").ClearFormats()

# G2V row (row 3): dSoH (ppm)
$ws.Range("G3").Value = "'50"
$ws.Range("G3").ClearFormats()

# V2G row (row 4): Trading rev.
$ws.Range("C4").Value = "'-0.9"
$ws.Range("C4").ClearFormats()

# V2G row (row 4): Adj. rev.
$ws.Range("F4").Value = "'-1.5"
$ws.Range("F4").ClearFormats()

# V2G row (row 4): dSoH (ppm)
$ws.Range("G4").Value = "'54.5"
$ws.Range("G4").ClearFormats()
